# Generate Report for Handoff
# Updates the localization-status report: marks rows 9-14 in the
# zh-cn and de-de sheets with Priority "ht", and refreshes the
# handoff/generate timestamps for the 983ef2cc-... file.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: Priority column (E) for rows 9-14 set to "ht"
for ($r = 9; $r -le 14; $r++) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
}

# de-de sheet: Priority column (E) for rows 9-14 set to "ht"
for ($r = 9; $r -le 14; $r++) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}

# zh-cn sheet: Latest Handoff Datetime column (H) for rows 9-14
for ($r = 9; $r -le 14; $r++) {
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-01 16:26:16"
}

# de-de sheet: Latest Handoff Datetime column (H) for rows 9-14
for ($r = 9; $r -le 14; $r++) {
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-01 16:26:21"
}

# Overview sheet: Latest HO Xliff Generate Date (G) for row 9
$wsOverview.Cells.Item(9, 7).Value = "2016-09-01 16:26:21"
